$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.6994255
$ws.Range("H2").Value = 1.398851
$ws.Range("I2").Value = 0.182573212173366
$ws.Range("J2").Value = 0.1466317745802101
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.0568385
$ws.Range("N2").Value = 0.113677
$ws.Range("O2").Value = 0.01274651757362603
$ws.Range("P2").Value = 0.008533937711420974
$ws.Range("Q2").Value = 0.03975429628175
$ws.Range("R2").Value = 0.159017185127
$ws.Range("S2").Value = 0.002327172657441164
$ws.Range("T2").Value = 0.001251346430782634
$ws.Range("G3").Value = 0.6994255
$ws.Range("H3").Value = 1.398851
$ws.Range("I3").Value = 0.182573212173366
$ws.Range("J3").Value = 0.1466317745802101
$ws.Range("O3").Value = 0.9872534824263741
$ws.Range("P3").Value = 0.9914660622885791
$ws.Range("Q3").Value = 3.0790815780755
$ws.Range("R3").Value = 18.474489468453
$ws.Range("S3").Value = 0.1802460395159249
$ws.Range("T3").Value = 0.1453804281494274
$ws.Range("I4").Value = 0.09857314952542007
$ws.Range("J4").Value = 0.1187519981886584
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.0568385
$ws.Range("N4").Value = 0.113677
$ws.Range("O4").Value = 0.01274651757362603
$ws.Range("P4").Value = 0.008533937711420974
$ws.Range("Q4").Value = 0.0214637522395
$ws.Range("R4").Value = 0.128782513437
$ws.Range("S4").Value = 0.001256464382713434
$ws.Range("T4").Value = 0.001013422155648787
$ws.Range("I5").Value = 0.09857314952542007
$ws.Range("J5").Value = 0.1187519981886584
$ws.Range("O5").Value = 0.9872534824263741
$ws.Range("P5").Value = 0.9914660622885791
$ws.Range("S5").Value = 0.09731668514270665
$ws.Range("T5").Value = 0.1177385760330096
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2206816666666667
$ws.Range("H6").Value = 0.662045
$ws.Range("I6").Value = 0.05760522135825098
$ws.Range("J6").Value = 0.06939755070551128
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.5
$ws.Range("M6").Value = 0.0568385
$ws.Range("N6").Value = 0.113677
$ws.Range("O6").Value = 0.01274651757362603
$ws.Range("P6").Value = 0.008533937711420974
$ws.Range("Q6").Value = 0.01254321491083333
$ws.Range("R6").Value = 0.075259289465
$ws.Range("S6").Value = 0.0007342659663755638
$ws.Range("T6").Value = 0.0005922343750460119
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.2206816666666667
$ws.Range("H7").Value = 0.662045
$ws.Range("I7").Value = 0.05760522135825098
$ws.Range("J7").Value = 0.06939755070551128
$ws.Range("O7").Value = 0.9872534824263741
$ws.Range("P7").Value = 0.9914660622885791
$ws.Range("Q7").Value = 0.9715071218483334
$ws.Range("R7").Value = 8.743564096635
$ws.Range("S7").Value = 0.05687095539187542
$ws.Range("T7").Value = 0.06880531633046527
$ws.Range("G8").Value = 1.2534795
$ws.Range("H8").Value = 2.506959
$ws.Range("I8").Value = 0.327199649867591
$ws.Range("J8").Value = 0.2627869923028463
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.5
$ws.Range("M8").Value = 0.0568385
$ws.Range("N8").Value = 0.113677
$ws.Range("O8").Value = 0.01274651757362603
$ws.Range("P8").Value = 0.008533937711420974
$ws.Range("Q8").Value = 0.07124589456075001
$ws.Range("R8").Value = 0.284983578243
$ws.Range("S8").Value = 0.004170656087121533
$ws.Range("T8").Value = 0.002242607823684153
$ws.Range("G9").Value = 1.2534795
$ws.Range("H9").Value = 2.506959
$ws.Range("I9").Value = 0.327199649867591
$ws.Range("J9").Value = 0.2627869923028463
$ws.Range("O9").Value = 0.9872534824263741
$ws.Range("P9").Value = 0.9914660622885791
$ws.Range("Q9").Value = 5.518194056329501
$ws.Range("R9").Value = 33.10916433797701
$ws.Range("S9").Value = 0.3230289937804695
$ws.Range("T9").Value = 0.2605443844791622
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.2558096666666667
$ws.Range("H10").Value = 0.767429
$ws.Range("I10").Value = 0.06677479238079162
$ws.Range("J10").Value = 0.08044421895849953
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.5
$ws.Range("M10").Value = 0.0568385
$ws.Range("N10").Value = 0.113677
$ws.Range("O10").Value = 0.01274651757362603
$ws.Range("P10").Value = 0.008533937711420974
$ws.Range("Q10").Value = 0.01453983773883333
$ws.Range("R10").Value = 0.08723902643300001
$ws.Range("S10").Value = 0.0008511460645569901
$ws.Range("T10").Value = 0.0006865059538357452
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.2558096666666667
$ws.Range("H11").Value = 0.767429
$ws.Range("I11").Value = 0.06677479238079162
$ws.Range("J11").Value = 0.08044421895849953
$ws.Range("O11").Value = 0.9872534824263741
$ws.Range("P11").Value = 0.9914660622885791
$ws.Range("Q11").Value = 1.126151151376333
$ws.Range("R11").Value = 10.135360362387
$ws.Range("S11").Value = 0.06592364631623464
$ws.Range("T11").Value = 0.07975771300466379
$ws.Range("G12").Value = 1.023908333333333
$ws.Range("H12").Value = 3.071725
$ws.Range("I12").Value = 0.2672739746945804
$ws.Range("J12").Value = 0.3219874652642745
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.5
$ws.Range("M12").Value = 0.0568385
$ws.Range("N12").Value = 0.113677
$ws.Range("O12").Value = 0.01274651757362603
$ws.Range("P12").Value = 0.008533937711420974
$ws.Range("Q12").Value = 0.05819741380416666
$ws.Range("R12").Value = 0.349184482825
$ws.Range("S12").Value = 0.003406812415417348
$ws.Range("T12").Value = 0.002747820972423644
$ws.Range("G13").Value = 1.023908333333333
$ws.Range("H13").Value = 3.071725
$ws.Range("I13").Value = 0.2672739746945804
$ws.Range("J13").Value = 0.3219874652642745
$ws.Range("O13").Value = 0.9872534824263741
$ws.Range("P13").Value = 0.9914660622885791
$ws.Range("Q13").Value = 4.507552679741667
$ws.Range("R13").Value = 40.567974117675
$ws.Range("S13").Value = 0.2638671622791631
$ws.Range("T13").Value = 0.3192396442918509
